$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Phase 1: resize the BOM table to cover the 3 new rows ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A2:F29"))

# --- Phase 2: write all cells that reuse already-existing text/values ---
# (order here does not affect the shared-string table, so rows are written top to bottom)
# Row 3
$ws.Range('C3').Value = 2
$ws.Range('D3').Value = 'VCC'
$ws.Range('E3').Value = 'LFB050CTP'
$ws.Range('F3').Value = 'Front face center LED light pipe'

# Row 4
$ws.Range('A4').Value = '1052620001'
$ws.Range('B4').Value = 'Antenna, 866MHz, 915MHz  Flat Patch, 0.4dBi, 1.4dBi U.FL (UMCC), IPEX MHF1 Adhesive'
$ws.Range('C4').Value = 1
$ws.Range('D4').Value = 'Molex'
$ws.Range('E4').Value = '1052620001'

# Row 5
$ws.Range('A5').Value = 'product-label'
$ws.Range('B5').Value = 'Fasson Transfer PET 75 / S8015 Matt silver heavy duty label 40mm x 25 mm'
$ws.Range('C5').Value = 1
$ws.Range('D5').Value = 'Fasson'
$ws.Range('E5').Value = 'PET75/S8015'
$ws.Range('F5').Value = 'Silver label applied on back of the bottom plastic case'

# Row 6
$ws.Range('A6').Value = 'Vayu-front_shell'
$ws.Range('B6').Value = 'Mec, Front shell'
$ws.Range('C6').Value = 1
$ws.Range('D6').Value = '3d print'
$ws.Range('E6').Value = 'Vayu-front_shell'
$ws.Range('F6').ClearContents()

# Row 7
$ws.Range('A7').Value = 'Vayu-wall_mounting_frame'
$ws.Range('B7').Value = 'Mec, Vayu wall mounting frame'
$ws.Range('C7').Value = 1
$ws.Range('D7').Value = '3d print'
$ws.Range('F7').ClearContents()

# Row 8
$ws.Range('A8').Value = 'Vayu-fan-shroude'
$ws.Range('B8').Value = 'Mec, Vayu fan shroude'
$ws.Range('C8').Value = 1
$ws.Range('D8').Value = '3d print'
$ws.Range('F8').ClearContents()

# Row 9
$ws.Range('A9').Value = 'Vayu-filter_cap_front'
$ws.Range('B9').Value = 'Mec, Vayu filter cap front'
$ws.Range('C9').Value = 1
$ws.Range('D9').Value = '3d print'
$ws.Range('F9').ClearContents()

# Row 10
$ws.Range('A10').Value = 'heat_accumulator_filter_cap_back'
$ws.Range('B10').Value = 'Mec, heat-accumulator filter cap back'
$ws.Range('C10').Value = 1
$ws.Range('D10').Value = '3d print'
$ws.Range('F10').ClearContents()

# Row 11
$ws.Range('A11').Value = 'heat_accumulator_shroud'
$ws.Range('B11').Value = 'Mec, heat-accumulator shroud'
$ws.Range('C11').Value = 2
$ws.Range('D11').Value = '3d print'
$ws.Range('F11').ClearContents()

# Row 12
$ws.Range('A12').Value = 'Vayu-filter'
$ws.Range('B12').Value = 'Filter pad, 96 mm diameter'
$ws.Range('C12').Value = 2
$ws.Range('D12').Value = 'Flexit'
$ws.Range('E12').Value = 'friskluftsventil 100FF'
$ws.Range('F12').Value = 'Custom cut to 96mm diameter'

# Row 13
$ws.Range('A13').Value = 'PVC-pipe'
$ws.Range('B13').Value = 'Pipe, 100m diameter. 300mm length'
$ws.Range('C13').Value = 1
$ws.Range('D13').Value = 'Flexit'
$ws.Range('E13').Value = 'Flexit RG100'
$ws.Range('F13').Value = 'Custom cut to wall length'

# Row 14
$ws.Range('A14').Value = 'M4 screw fan'
$ws.Range('B14').Value = 'screw,M4 x 16mm, self tapping plastic screw'
$ws.Range('C14').Value = 2
$ws.Range('D14').Value = 'Amazon'
$ws.Range('E14').ClearContents()
$ws.Range('F14').ClearContents()

# Row 15
$ws.Range('A15').Value = 'Anchor'
$ws.Range('B15').Value = 'Super plug, NY, 5mm hole, 25mm length, EXPANDET'
$ws.Range('C15').Value = 4
$ws.Range('D15').Value = 'Amazon'
$ws.Range('E15').ClearContents()
$ws.Range('F15').ClearContents()

# Row 16
$ws.Range('A16').Value = 'M4 wall screw'
$ws.Range('B16').Value = 'screw, M4 x 16mm length'
$ws.Range('C16').Value = 4
$ws.Range('D16').Value = 'Amazon'
$ws.Range('E16').ClearContents()
$ws.Range('F16').ClearContents()

# Row 17
$ws.Range('A17').Value = 'M2 screw'
$ws.Range('B17').Value = 'screw, M2.0, 5mm Length, Pan head, Philips OR Rouded head Torx Thread-forming screw for plastic'
$ws.Range('C17').Value = 8
$ws.Range('D17').Value = 'Amazon'
$ws.Range('E17').ClearContents()
$ws.Range('F17').ClearContents()

# Row 18
$ws.Range('A18').Value = 'M2 Nut'
$ws.Range('B18').Value = 'Nut, M2, Hex nut, 1.5mm thick, 3.8mm length'
$ws.Range('C18').Value = 4
$ws.Range('D18').Value = 'Amazon'
$ws.Range('E18').ClearContents()
$ws.Range('F18').ClearContents()

# Row 19
$ws.Range('A19').Value = 'SanAce_9RF_92x38_RBD'
$ws.Range('B19').Value = 'Fan Tubeaxial 12VDC Round - 92mm Dia Ball 42.4 CFM (1.19m³/min) 4 Wire Leads'
$ws.Range('C19').Value = 1
$ws.Range('D19').Value = 'Sanyo Denki America Inc.'
$ws.Range('E19').Value = '9RF0912P1H001'
$ws.Range('F19').ClearContents()

# Row 20
$ws.Range('A20').Value = 'Vayu-R1A-PCBA'
$ws.Range('C20').Value = 1
$ws.Range('D20').Value = 'Self assembled'
$ws.Range('E20').Value = 'Vayu-R1A-PCBA'
$ws.Range('F20').ClearContents()

# Row 21
$ws.Range('A21').Value = 'Vayu_i2c-PCBA'
$ws.Range('C21').Value = 1
$ws.Range('D21').Value = 'Self assembled'
$ws.Range('F21').ClearContents()

# Row 22
$ws.Range('A22').Value = 'Vayu_AC_DC-PCBA'
$ws.Range('C22').Value = 1
$ws.Range('D22').Value = 'Self assembled'
$ws.Range('F22').ClearContents()

# Row 23
$ws.Range('B23').Value = 'Tube, Aluminim, ID4mm, OD6mm, 166mm length'
$ws.Range('C23').Value = 211
$ws.Range('D23').Value = 'Kaiserthal'
$ws.Range('E23').Value = 'N/A'
$ws.Range('F23').Value = 'Bought from Hornbach and custom cut to 166mm'

# Row 24
$ws.Range('A24').Value = 'JST cable connector'
$ws.Range('B24').Value = '4 Position Cable Assembly Rectangular Socket to Socket, Reversed 0.667'' (203.20mm, 8.00")'
$ws.Range('C24').Value = 1
$ws.Range('D24').Value = 'JST Sales America Inc.'
$ws.Range('E24').Value = 'A04SR04SR30K203A'
$ws.Range('F24').ClearContents()

# Row 25
$ws.Range('A25').Value = 'AC-plug'
$ws.Range('B25').Value = 'plug, 2 pin AC cord, EU type, 2.5 meter'
$ws.Range('C25').Value = 1
$ws.Range('D25').Value = 'Amazon'
$ws.Range('E25').ClearContents()
$ws.Range('F25').ClearContents()

# Row 26
$ws.Range('A26').Value = 'magnet'
$ws.Range('B26').Value = 'magnet, Neodeam, 20mm x 10mm x 2mm'
$ws.Range('C26').Value = 1
$ws.Range('D26').Value = 'Amazon'
$ws.Range('E26').ClearContents()
$ws.Range('F26').ClearContents()

# Row 27
$ws.Range('C27').Value = 1
$ws.Range('F27').ClearContents()

# Row 28
$ws.Range('C28').Value = 1
$ws.Range('D28').Value = 'Amazon'
$ws.Range('E28').ClearContents()

# Row 29
$ws.Range('C29').Value = 1
$ws.Range('D29').Value = 'Amazon'
$ws.Range('E29').ClearContents()
$ws.Range('F29').ClearContents()

# --- Phase 3: write the cells that introduce brand-new text, in the exact order
#     the values were first typed (so the shared-string table order matches) ---
$ws.Range('A27').Value = 'Aluminium tape'
$ws.Range('A23').Value = 'Aluminium tubes'
$ws.Range('B27').Value = 'tape, 50mm width, 250mm length, Metallic searing tape for high and low temperature channels, waterproof'
$ws.Range('D27').Value = 'Tesa'
$ws.Range('E27').Value = 'Tesa-60672'
$ws.Range('A28').Value = 'NTC-100K-probe'
$ws.Range('B28').Value = 'probe, NTC temperature probe, 100K, B3950 Thermistor, OD3mm Cartridge sensor'
$ws.Range('F28').Value = 'Used in 3d printer extruder hot ends'
$ws.Range('A29').Value = 'Polymide-film-heater'
$ws.Range('B29').Value = 'Heater, Polymide film, 12Volt, 10watts, 30mm x 90mm'
$ws.Range('A3').Value = 'LFB075CTP'
$ws.Range('B3').Value = 'Light Pipe, 3mm dia, 19mm length, Single Clear, Diffused Rigid Panel Mount, Press Fit, Front'
$ws.Range('F4').Value = 'Optinal Sub-GHz antenna'
$ws.Range('E21').Value = 'Vayu-i2c_R1A-PCBA'
$ws.Range('E22').Value = 'Vayu_AC_DC-R1A-PCBA'
$ws.Range('E7').Value = 'Vayu_wall_frame'
$ws.Range('E8').Value = 'Vayu_fan_shroude'
$ws.Range('E9').Value = 'Vayu_filter_cap_front'
$ws.Range('E11').Value = 'Vayu_heat_accumulator_shroude'
$ws.Range('E10').Value = 'Vayu_filter_cap_back'
$ws.Range('B20').Value = 'PCBA, Vayu-R1A, Assembled'
$ws.Range('B21').Value = 'PCBA, Vayu-i2c-R1A, Assembled'
$ws.Range('B22').Value = 'PCBA, Vayu_AC_DC-R1A, Assembled'

# --- Final touches ---
$ws.Range("A29").Select()

Write-Host "done"
